$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shape = $s.Shapes.Item(4)
$tr = $shape.TextFrame.TextRange

# Remove the empty paragraph (endParaRPr sz=2000, hlinkClick rId7) that sits
# right after the "https://github.com/StabiBerlin/sbb-relevance-test" line.
# Deleting the paragraph that currently occupies index 4 merges its paragraph
# mark forward so that the surviving paragraph (previously index 3, with the
# "sz=2000 dirty=0" end run properties) shifts up to become the new index 3,
# exactly matching the target edit, while leaving the URL paragraph untouched.
$tr.Paragraphs(4).Delete()
